$d = $word.ActiveDocument

# ===========================================================================
# Helper functions
# ===========================================================================

function Replace-ParagraphText($para, $newText) {
    # Rewrites a (single- or multi-run) paragraph's text in one run,
    # preserving the paragraph mark itself (and anything tied to it).
    $r = $para.Range
    $r.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; exclude the pilcrow
    $r.Delete()
    $r.InsertAfter($newText)
}

function Set-ParagraphTextKeepingTrailer($para, $newText) {
    # Replaces paragraph text but keeps any trailing elements (e.g.
    # bookmarks) anchored after the last run (used when a bookmark sits at
    # the paragraph end, right before the pilcrow).
    $r = $para.Range
    $r.MoveEnd(1, -1) | Out-Null
    $origStart = $r.Start
    $origEnd = $r.End
    $r.InsertBefore($newText)
    $oldStart = $origStart + $newText.Length
    $oldEnd = $origEnd + $newText.Length
    $d.Range($oldStart, $oldEnd).Delete()
}

function New-ParagraphAfter($paraIndex, $newText) {
    # Inserts a new plain paragraph (single run) after paragraph number
    # $paraIndex (1-based). Returns the 1-based index of the new paragraph.
    $para = $d.Paragraphs.Item($paraIndex)
    $r = $para.Range
    $r.InsertParagraphAfter()
    $newIndex = $paraIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $nr = $newPara.Range
    $nr.InsertAfter($newText)
    return $newIndex
}

Write-Host "Paragraphs at start: $($d.Paragraphs.Count)"

# ===========================================================================
# STEP A: remove the old INTRO section (4 paragraphs): "INTRO", "Interstitial
# Cystitis is a chronic condition...", "The easiest way to think...",
# "Unfortunately at this stage...".
# ===========================================================================
$d.Content.Find.Execute("INTRO^p", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$fr = $d.Content
$fr.Find.Execute("Interstitial Cystitis is a chronic", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $fr.Start

$fr2 = $d.Content
$fr2.Find.Execute("infection or allergies.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $fr2.End + 1  # swallow the trailing paragraph mark too

$d.Range($startPos, $endPos).Delete()

Write-Host "After step A: $($d.Paragraphs.Count) paragraphs"

# ===========================================================================
# STEP B: insert the first interview answer after "Firstly can you explain...".
# ===========================================================================
$answer1 = 'Interstitial Cystitis (IC) is a chronic painful bladder condition. I was finally diagnosed with it last year (2022) after years of going back to doctor after doctor and surgery after surgery for endometriosis and no one would understand why I was still having so much pain in my daily life. I had many doctors tell me it was all in my head until finally I found a doctor that would listen to me and went deeper to find the source. Since finding out I now suffer from 2 chronic illness I have had to change certain things in my life. I have changed my diet in ways that I try to stay away or at least limit my intake of what triggers my flare ups. However it’s not as easy as that sounds unfortunately. I have learnt to listen to my body and react accordingly. Rest when I use to push myself. I now drink water like it’s going out of fashion as that helps to flush me out so when I do eat and drink certain foods the acid from them doesn’t linger in my bladder.'

$idxFirstly = 3
New-ParagraphAfter $idxFirstly $answer1 | Out-Null
Write-Host "After step B: $($d.Paragraphs.Count) paragraphs"

# ===========================================================================
# STEP C: insert the second interview answer after "How would this app help...".
# ===========================================================================
$answer2 = 'This app would help me immensely as my biggest struggle I have found is remembering to document my intake. Pervious apps I have tried I end up giving up on as we all know life gets in the way and I found sometimes I could go a day or so before I got a chance to get back to it and then I was racking my brain to remember what I had eaten and was the pain at this time or that time. The other struggle I have is trying to work out what actually causes flare ups as I soon found out with trying to eliminate foods is I would eat something and be fine one day and think okay this is good than a few days to a week later I would eat the same food and have a completely different reaction to it so I definitely believe that certain times of the month affect me more so than others. I have been taking medication that my doctor prescribed me and that has seemed to help a lot as well as the rule of eat small amount of trigger foods and then immediately after drink 1-2 cups of water.   '

$idxHowWould = 5
New-ParagraphAfter $idxHowWould $answer2 | Out-Null
Write-Host "After step C: $($d.Paragraphs.Count) paragraphs"

# ===========================================================================
# STEP D: the old "Bladder instillations are administered..." paragraph
# becomes the new "The easy and non-invasive treatment plans..." paragraph.
# ===========================================================================
$bladderPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.StartsWith("Bladder instillations are administered")) {
        $bladderPara = $pp
        $bladderIdx = $i
        break
    }
}
Replace-ParagraphText $bladderPara 'The easy and non-invasive treatment plans to reduce the risk or flare ups include maintaining hydration.'
Write-Host "After step D: $($d.Paragraphs.Count) paragraphs, idx=$bladderIdx text=[$($d.Paragraphs.Item($bladderIdx).Range.Text)]"

# ===========================================================================
# STEP E: insert three new paragraphs right after it.
# ===========================================================================
$idx = $bladderIdx
$idx = New-ParagraphAfter $idx 'Drinking plenty of water can water down the acidity of most food groups which can ease the irritation and pain.'
$idx = New-ParagraphAfter $idx 'Staying away from acidic food groups. The acid from the food can cause pain against the bladder walls giving burning sensations. Carbonated drinks, alcohol, caffeine, citrus products and foods containing high concentrations of vitamin C can all cause issues.'
$idx = New-ParagraphAfter $idx 'Reducing stress in your life is a good tip for anyone but especially IC sufferers. Smoking should be stopped as it can also increase pain and exercise regularly.'
Write-Host "After step E: $($d.Paragraphs.Count) paragraphs, last idx=$idx"
